# Remove the "Compartment" column from the Submodels worksheet
# (wc_lang.Submodel no longer has a compartment attribute).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Submodels")
$ws.Activate()

# Column D is "Compartment" - delete the entire column, shifting the
# remaining columns (Biomass reaction, Objective function, Comments,
# References) left.
$ws.Range("D1").EntireColumn.Delete()
